$d = $word.ActiveDocument

# Locate the exact text " e crescer ;" (end of the "Está na empresa de
# grande prestígio e crescer ;" objective bullet).
$m = $d.Content
$found = $m.Find.Execute(" e crescer ;")
if (-not $found) {
    throw "Could not locate ' e crescer ;'"
}
$start = $m.Start
$end = $m.End

# " e crescer ;" -> split after " e " (space, e, space) so the bookmark
# lands between " e " and "crescer ;".
$splitPoint = $start + 3

# Relocate the "_GoBack" bookmark to this split point. Bookmark names
# are unique in Word, so re-adding "_GoBack" here automatically removes
# it from its old location (near "e sempre aprender") and creates it
# here, splitting the run in two around the (now collapsed) bookmark.
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Tighten "crescer ;" to "crescer;" by deleting the space before the
# semicolon (last two characters of the original match).
$spaceRange = $d.Range($end - 2, $end - 1)
$spaceRange.Text = ""
